# Update CodeSystem-PueblosOriginariosCS metadata: new version (status, experimental flag, date)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> false (leading apostrophe forces plain text, not a boolean;
# re-paste the original cell formatting so the style index is unaffected)
$ws.Range("B7").Value = "'false"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Date: refreshed publish timestamp
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"
